$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 33998
$ws.Range("E2").Value = 2734
$ws.Range("F2").Value = 2734
$ws.Range("G2").Value = 4751
$ws.Range("H2").Value = 3267
$ws.Range("I2").Value = 3257
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 86545
$ws.Range("L2").Value = 24252
$ws.Range("M2").Value = 62293
$ws.Range("N2").Value = 61970
$ws.Range("O2").Value = 324
$ws.Range("P2").Value = 563
$ws.Range("Q2").Value = 2765
$ws.Range("R2").Value = -1648
$ws.Range("S2").Value = -599
$ws.Range("T2").Value = 3206
$ws.Range("U2").Value = -441
$ws.Range("V2").Value = 10248
$ws.Range("W2").Value = 8.039999999999999
$ws.Range("X2").Value = 9.609999999999999
$ws.Range("Y2").Value = 5.79
$ws.Range("Z2").Value = 4.16
$ws.Range("AA2").Value = 38.93
$ws.Range("AB2").Value = 8390.58
$ws.Range("AC2").Value = 30964
$ws.Range("AD2").Value = 16.79
$ws.Range("AE2").Value = 632540
$ws.Range("AF2").Value = 0.82
$ws.Range("AG2").Value = 9000
$ws.Range("AH2").Value = 1.73
$ws.Range("AI2").Value = 27.07
$ws.Range("AJ2").Value = 10520000

# Row 3
$ws.Range("D3").Value = 34144
$ws.Range("E3").Value = 3092
$ws.Range("F3").Value = 3092
$ws.Range("G3").Value = 2507
$ws.Range("H3").Value = 1851
$ws.Range("I3").Value = 1836
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 86684
$ws.Range("L3").Value = 27998
$ws.Range("M3").Value = 58687
$ws.Range("N3").Value = 58579
$ws.Range("O3").Value = 107
$ws.Range("P3").Value = 564
$ws.Range("Q3").Value = 954
$ws.Range("R3").Value = -8647
$ws.Range("S3").Value = 4639
$ws.Range("T3").Value = 3603
$ws.Range("U3").Value = -2649
$ws.Range("V3").Value = 16450
$ws.Range("W3").Value = 9.06
$ws.Range("X3").Value = 5.42
$ws.Range("Y3").Value = 3.04
$ws.Range("Z3").Value = 2.14
$ws.Range("AA3").Value = 47.71
$ws.Range("AB3").Value = 8559.48
$ws.Range("AC3").Value = 17443
$ws.Range("AD3").Value = 23.93
$ws.Range("AE3").Value = 595712
$ws.Range("AF3").Value = 0.7
$ws.Range("AG3").Value = 9000
$ws.Range("AH3").Value = 2.16
$ws.Range("AI3").Value = 48.19
$ws.Range("AJ3").Value = 10556513

# Row 4
$ws.Range("D4").Value = 34905
$ws.Range("E4").Value = 3266
$ws.Range("F4").Value = 3266
$ws.Range("G4").Value = 2320
$ws.Range("H4").Value = 1530
$ws.Range("I4").Value = 1523
$ws.Range("J4").Value = 6
$ws.Range("K4").Value = 91624
$ws.Range("L4").Value = 31954
$ws.Range("M4").Value = 59669
$ws.Range("N4").Value = 59560
$ws.Range("O4").Value = 110
$ws.Range("P4").Value = 564
$ws.Range("Q4").Value = 4311
$ws.Range("R4").Value = -4980
$ws.Range("S4").Value = 2019
$ws.Range("T4").Value = 2905
$ws.Range("U4").Value = 1406
$ws.Range("V4").Value = 20005
$ws.Range("W4").Value = 9.359999999999999
$ws.Range("X4").Value = 4.38
$ws.Range("Y4").Value = 2.58
$ws.Range("Z4").Value = 1.72
$ws.Range("AA4").Value = 53.55
$ws.Range("AB4").Value = 8696.07
$ws.Range("AC4").Value = 14430
$ws.Range("AD4").Value = 24.91
$ws.Range("AE4").Value = 605683
$ws.Range("AF4").Value = 0.59
$ws.Range("AG4").Value = 9000
$ws.Range("AH4").Value = 2.5
$ws.Range("AI4").Value = 58.1
$ws.Range("AJ4").Value = 10556513

# Row 5
$ws.Range("D5").Value = 38640
$ws.Range("E5").Value = 3298
$ws.Range("F5").Value = 3298
$ws.Range("G5").Value = 939
$ws.Range("H5").Value = 423
$ws.Range("I5").Value = 418
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 95438
$ws.Range("L5").Value = 36775
$ws.Range("M5").Value = 58663
$ws.Range("N5").Value = 58552
$ws.Range("O5").Value = 110
$ws.Range("P5").Value = 564
$ws.Range("Q5").Value = 3329
$ws.Range("R5").Value = -6585
$ws.Range("S5").Value = 1416
$ws.Range("T5").Value = 4223
$ws.Range("U5").Value = -894
$ws.Range("V5").Value = 22486
$ws.Range("W5").Value = 8.539999999999999
$ws.Range("X5").Value = 1.09
$ws.Range("Y5").Value = 0.71
$ws.Range("Z5").Value = 0.45
$ws.Range("AA5").Value = 62.69
$ws.Range("AB5").Value = 10675.74
$ws.Range("AC5").Value = 3958
$ws.Range("AD5").Value = 96.26000000000001
$ws.Range("AE5").Value = 595438
$ws.Range("AF5").Value = 0.64
$ws.Range("AG5").Value = 9000
$ws.Range("AH5").Value = 2.36
$ws.Range("AI5").Value = 211.8
$ws.Range("AJ5").Value = 10556513

# Row 6
$ws.Range("D6").Value = 37822
$ws.Range("E6").Value = 2435
$ws.Range("F6").Value = 2435
$ws.Range("G6").Value = -123
$ws.Range("H6").Value = -231
$ws.Range("I6").Value = -236
$ws.Range("K6").Value = 89654
$ws.Range("L6").Value = 32241
$ws.Range("M6").Value = 57413
$ws.Range("N6").Value = 57300
$ws.Range("P6").Value = 564
$ws.Range("Q6").Value = 5390
$ws.Range("R6").Value = -503
$ws.Range("S6").Value = -5238
$ws.Range("T6").Value = 2360
$ws.Range("U6").Value = 3030
$ws.Range("V6").Value = 18757
$ws.Range("W6").Value = 6.44
$ws.Range("X6").Value = -0.61
$ws.Range("Y6").Value = -0.41
$ws.Range("Z6").Value = -0.25
$ws.Range("AA6").Value = 56.16
$ws.Range("AB6").Value = 10476.26
$ws.Range("AC6").Value = -2233
$ws.Range("AD6").Value = -138.37
$ws.Range("AE6").Value = 582704
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 9000
$ws.Range("AH6").Value = 2.91
$ws.Range("AI6").Value = -375.41
$ws.Range("AJ6").Value = 10556513

# Row 7
$ws.Range("D7").Value = 33330
$ws.Range("E7").Value = 1671
$ws.Range("G7").Value = -742
$ws.Range("H7").Value = -748
$ws.Range("I7").Value = -766
$ws.Range("K7").Value = 91924
$ws.Range("L7").Value = 36004
$ws.Range("M7").Value = 55920
$ws.Range("N7").Value = 55254
$ws.Range("P7").Value = 542
$ws.Range("Q7").Value = 2670
$ws.Range("R7").Value = -6483
$ws.Range("S7").Value = 3634
$ws.Range("T7").Value = 2394
$ws.Range("U7").Value = -1188
$ws.Range("W7").Value = 5.01
$ws.Range("X7").Value = -2.24
$ws.Range("Y7").Value = -1.36
$ws.Range("Z7").Value = -0.82
$ws.Range("AA7").Value = 64.39
$ws.Range("AC7").Value = -7251
$ws.Range("AD7").Value = -27.79
$ws.Range("AE7").Value = 561902
$ws.Range("AF7").Value = 0.36
$ws.Range("AG7").Value = 9000
$ws.Range("AH7").Value = 4.47
$ws.Range("AI7").Value = -104.48

# Row 8
$ws.Range("D8").Value = 44970
$ws.Range("E8").Value = 2994
$ws.Range("G8").Value = 2334
$ws.Range("H8").Value = 1738
$ws.Range("I8").Value = 1374
$ws.Range("K8").Value = 103800
$ws.Range("L8").Value = 47048
$ws.Range("M8").Value = 56753
$ws.Range("N8").Value = 55701
$ws.Range("P8").Value = 522
$ws.Range("Q8").Value = 1110
$ws.Range("R8").Value = -4710
$ws.Range("S8").Value = -1620
$ws.Range("T8").Value = 4114
$ws.Range("U8").Value = -7472
$ws.Range("W8").Value = 6.66
$ws.Range("X8").Value = 3.86
$ws.Range("Y8").Value = 2.48
$ws.Range("Z8").Value = 1.78
$ws.Range("AA8").Value = 82.90000000000001
$ws.Range("AC8").Value = 15457
$ws.Range("AD8").Value = 13.04
$ws.Range("AE8").Value = 682325
$ws.Range("AF8").Value = 0.3
$ws.Range("AG8").Value = 9000
$ws.Range("AH8").Value = 4.47
$ws.Range("AI8").Value = 58.2

# Row 9
$ws.Range("D9").Value = 46783
$ws.Range("E9").Value = 3286
$ws.Range("G9").Value = 2833
$ws.Range("H9").Value = 2108
$ws.Range("I9").Value = 1663
$ws.Range("K9").Value = 104474
$ws.Range("L9").Value = 46739
$ws.Range("M9").Value = 57735
$ws.Range("N9").Value = 56163
$ws.Range("P9").Value = 522
$ws.Range("Q9").Value = 4899
$ws.Range("R9").Value = -1835
$ws.Range("S9").Value = -1955
$ws.Range("T9").Value = 2939
$ws.Range("U9").Value = 1979
$ws.Range("W9").Value = 7.02
$ws.Range("X9").Value = 4.51
$ws.Range("Y9").Value = 2.97
$ws.Range("Z9").Value = 2.02
$ws.Range("AA9").Value = 80.95
$ws.Range("AC9").Value = 18711
$ws.Range("AD9").Value = 10.77
$ws.Range("AE9").Value = 687981
$ws.Range("AF9").Value = 0.29
$ws.Range("AG9").Value = 9000
$ws.Range("AH9").Value = 4.47
$ws.Range("AI9").Value = 48.1
